# Team05Report 10.23 update — adds user stories US27 (Correct gender of family
# members) and US28 (Delete bigamous-parent records) across the Backlog,
# Burndown, Sprint3 and Stories sheets.

$wb = $excel.ActiveWorkbook

$usCorrectGenderShort = "Correct the gender of the family members "
$usCorrectGenderLong  = "Correct the wrong gender of family members and display the right"
$usBigamyShort        = "Delete the information of parents who are bigamy at the same time"
$usBigamyLong         = "Delete the data that parents who are bigamy at the same time"

# ---------------------------------------------------------------------
# Backlog: append the two new backlog items as rows 28-29
# ---------------------------------------------------------------------
$backlog = $wb.Worksheets.Item("Backlog")

$backlog.Cells.Item(28, 1).Value = 3
$backlog.Cells.Item(28, 2).Value = "US27"
$backlog.Cells.Item(28, 3).Value = $usCorrectGenderShort
$backlog.Cells.Item(28, 4).Value = "Rx"
$backlog.Cells.Item(28, 5).Value = "Coding"

$backlog.Cells.Item(29, 1).Value = 3
$backlog.Cells.Item(29, 2).Value = "US28"
$backlog.Cells.Item(29, 3).Value = $usBigamyShort
$backlog.Cells.Item(29, 4).Value = "Rx"
$backlog.Cells.Item(29, 5).Value = "Coding"

# ---------------------------------------------------------------------
# Burndown: record the 10/23 sprint checkpoint in row 4
# ---------------------------------------------------------------------
$burndown = $wb.Worksheets.Item("Burndown")

$burndown.Cells.Item(4, 1).Value = 42299
$burndown.Cells.Item(4, 1).NumberFormat = "m/d"
$burndown.Cells.Item(4, 2).Value = 20
$burndown.Cells.Item(4, 3).Value = 10
$burndown.Cells.Item(4, 4).Value = 240
$burndown.Cells.Item(4, 5).Value = 120
$burndown.Cells.Item(4, 6).Value = 120

# ---------------------------------------------------------------------
# Sprint3: fill in the story name / owner / status / estimates for the
# already-present US27 and US28 rows (8 and 9)
# ---------------------------------------------------------------------
$sprint3 = $wb.Worksheets.Item("Sprint3")

$sprint3.Cells.Item(8, 2).Value = $usCorrectGenderShort
$sprint3.Cells.Item(8, 3).Value = "Rx"
$sprint3.Cells.Item(8, 4).Value = "Coding"
$sprint3.Cells.Item(8, 5).Value = 40
$sprint3.Cells.Item(8, 6).Value = 75

$sprint3.Cells.Item(9, 2).Value = $usBigamyShort
$sprint3.Cells.Item(9, 3).Value = "Rx"
$sprint3.Cells.Item(9, 4).Value = "Coding"
$sprint3.Cells.Item(9, 5).Value = 40
$sprint3.Cells.Item(9, 6).Value = 75

# ---------------------------------------------------------------------
# Stories: document the two stories in rows 28-29
# ---------------------------------------------------------------------
$stories = $wb.Worksheets.Item("Stories")

$stories.Cells.Item(28, 1).Value = "US27"
$stories.Cells.Item(28, 2).Value = $usCorrectGenderShort
$stories.Cells.Item(28, 3).Value = $usCorrectGenderLong

$stories.Cells.Item(29, 1).Value = "US28"
$stories.Cells.Item(29, 2).Value = $usBigamyShort
$stories.Cells.Item(29, 3).Value = $usBigamyLong
